# Thailand Premier League workbook update
# - Swap the data of rows 15/16 (match order correction), keeping the "id" (col A) fixed per row.
# - Swap the data of rows 85/86 (match order correction), keeping the "id" (col A) fixed per row.
# - Append 5 new match rows (171-175) with odds/result data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Swap full row content (columns B..AC) between rows 15 and 16.
#    Column A (the sequential "id") stays put on its own row.
# ---------------------------------------------------------------------------
$row15 = $ws.Range("B15:AC15").Value2
$row16 = $ws.Range("B16:AC16").Value2
$ws.Range("B15:AC15").Value = $row16
$ws.Range("B16:AC16").Value = $row15

# ---------------------------------------------------------------------------
# 2) Swap full row content (columns B..AC) between rows 85 and 86.
# ---------------------------------------------------------------------------
$row85 = $ws.Range("B85:AC85").Value2
$row86 = $ws.Range("B86:AC86").Value2
$ws.Range("B85:AC85").Value = $row86
$ws.Range("B86:AC86").Value = $row85

# ---------------------------------------------------------------------------
# 3) Append new match rows 171-175.
#    First clone the formatting of the last existing row (170) down so the
#    new rows pick up the same styles (bold/bordered id column, date format
#    on the date column, etc.), then overwrite the values.
# ---------------------------------------------------------------------------
$lastRowFormat = $ws.Range("A170:AC170")
for ($r = 171; $r -le 175; $r++) {
    $lastRowFormat.Copy($ws.Range("A" + $r + ":AC" + $r))
}

# Row 171
$ws.Cells.Item(171, 1).Value = 169
$ws.Cells.Item(171, 2).Value = 6992693
$ws.Cells.Item(171, 3).Value = "Thailand Premier League"
$ws.Cells.Item(171, 4).Value = "Thailand Premier League"
$ws.Cells.Item(171, 5).Value = 45382.3125
$ws.Cells.Item(171, 6).Value = "Buriram United"
$ws.Cells.Item(171, 7).Value = "Chiangrai Utd"
$ws.Cells.Item(171, 8).Value = 2
$ws.Cells.Item(171, 9).Value = 1
$ws.Cells.Item(171, 10).Value = "H"
$ws.Cells.Item(171, 11).Value = 1.4
$ws.Cells.Item(171, 12).Value = 4.333
$ws.Cells.Item(171, 13).Value = 6
$ws.Cells.Item(171, 14).Value = 1.2
$ws.Cells.Item(171, 15).Value = 6
$ws.Cells.Item(171, 16).Value = 9
$ws.Cells.Item(171, 17).Value = -2
$ws.Cells.Item(171, 18).Value = 1.85
$ws.Cells.Item(171, 19).Value = 1.95
$ws.Cells.Item(171, 20).Value = 3.5
$ws.Cells.Item(171, 21).Value = 1.975
$ws.Cells.Item(171, 22).Value = 1.825
$ws.Cells.Item(171, 23).Value = 0.2
$ws.Cells.Item(171, 24).Value = -1
$ws.Cells.Item(171, 25).Value = -1
$ws.Cells.Item(171, 26).Value = -1
$ws.Cells.Item(171, 27).Value = 0.95
$ws.Cells.Item(171, 28).Value = -1
$ws.Cells.Item(171, 29).Value = 0.825

# Row 172
$ws.Cells.Item(172, 1).Value = 170
$ws.Cells.Item(172, 2).Value = 6992335
$ws.Cells.Item(172, 3).Value = "Thailand Premier League"
$ws.Cells.Item(172, 4).Value = "Thailand Premier League"
$ws.Cells.Item(172, 5).Value = 45382.33333333334
$ws.Cells.Item(172, 6).Value = "Muang Thong United"
$ws.Cells.Item(172, 7).Value = "Police Tero FC"
$ws.Cells.Item(172, 8).Value = 3
$ws.Cells.Item(172, 9).Value = 1
$ws.Cells.Item(172, 10).Value = "H"
$ws.Cells.Item(172, 11).Value = 1.444
$ws.Cells.Item(172, 12).Value = 4.2
$ws.Cells.Item(172, 13).Value = 5.5
$ws.Cells.Item(172, 14).Value = 1.363
$ws.Cells.Item(172, 15).Value = 4.5
$ws.Cells.Item(172, 16).Value = 6
$ws.Cells.Item(172, 17).Value = -1.5
$ws.Cells.Item(172, 18).Value = 1.975
$ws.Cells.Item(172, 19).Value = 1.825
$ws.Cells.Item(172, 20).Value = 3.25
$ws.Cells.Item(172, 21).Value = 1.875
$ws.Cells.Item(172, 22).Value = 1.925
$ws.Cells.Item(172, 23).Value = 0.363
$ws.Cells.Item(172, 24).Value = -1
$ws.Cells.Item(172, 25).Value = -1
$ws.Cells.Item(172, 26).Value = 0.9750000000000001
$ws.Cells.Item(172, 27).Value = -1
$ws.Cells.Item(172, 28).Value = 0.875
$ws.Cells.Item(172, 29).Value = -1

# Row 173
$ws.Cells.Item(173, 1).Value = 171
$ws.Cells.Item(173, 2).Value = 6992692
$ws.Cells.Item(173, 3).Value = "Thailand Premier League"
$ws.Cells.Item(173, 4).Value = "Thailand Premier League"
$ws.Cells.Item(173, 5).Value = 45382.35416666666
$ws.Cells.Item(173, 6).Value = "Bangkok United"
$ws.Cells.Item(173, 7).Value = "Chonburi"
$ws.Cells.Item(173, 8).Value = 6
$ws.Cells.Item(173, 9).Value = 0
$ws.Cells.Item(173, 10).Value = "H"
$ws.Cells.Item(173, 11).Value = 1.444
$ws.Cells.Item(173, 12).Value = 4.2
$ws.Cells.Item(173, 13).Value = 5.5
$ws.Cells.Item(173, 14).Value = 1.4
$ws.Cells.Item(173, 15).Value = 4.333
$ws.Cells.Item(173, 16).Value = 6
$ws.Cells.Item(173, 17).Value = -1.25
$ws.Cells.Item(173, 18).Value = 1.825
$ws.Cells.Item(173, 19).Value = 1.975
$ws.Cells.Item(173, 20).Value = 2.75
$ws.Cells.Item(173, 21).Value = 1.75
$ws.Cells.Item(173, 22).Value = 1.95
$ws.Cells.Item(173, 23).Value = 0.3999999999999999
$ws.Cells.Item(173, 24).Value = -1
$ws.Cells.Item(173, 25).Value = -1
$ws.Cells.Item(173, 26).Value = 0.825
$ws.Cells.Item(173, 27).Value = -1
$ws.Cells.Item(173, 28).Value = 0.75
$ws.Cells.Item(173, 29).Value = -1

# Row 174
$ws.Cells.Item(174, 1).Value = 172
$ws.Cells.Item(174, 2).Value = 6992688
$ws.Cells.Item(174, 3).Value = "Thailand Premier League"
$ws.Cells.Item(174, 4).Value = "Thailand Premier League"
$ws.Cells.Item(174, 5).Value = 45382.375
$ws.Cells.Item(174, 6).Value = "Khonkaen United"
$ws.Cells.Item(174, 7).Value = "Nakhon Pathom FC"
$ws.Cells.Item(174, 8).Value = 2
$ws.Cells.Item(174, 9).Value = 1
$ws.Cells.Item(174, 10).Value = "H"
$ws.Cells.Item(174, 11).Value = 2.7
$ws.Cells.Item(174, 12).Value = 3.6
$ws.Cells.Item(174, 13).Value = 2.15
$ws.Cells.Item(174, 14).Value = 2.1
$ws.Cells.Item(174, 15).Value = 3.5
$ws.Cells.Item(174, 16).Value = 2.875
$ws.Cells.Item(174, 17).Value = -0.25
$ws.Cells.Item(174, 18).Value = 1.9
$ws.Cells.Item(174, 19).Value = 1.9
$ws.Cells.Item(174, 20).Value = 2.75
$ws.Cells.Item(174, 21).Value = 1.825
$ws.Cells.Item(174, 22).Value = 1.975
$ws.Cells.Item(174, 23).Value = 1.1
$ws.Cells.Item(174, 24).Value = -1
$ws.Cells.Item(174, 25).Value = -1
$ws.Cells.Item(174, 26).Value = 0.8999999999999999
$ws.Cells.Item(174, 27).Value = -1
$ws.Cells.Item(174, 28).Value = 0.4125
$ws.Cells.Item(174, 29).Value = -0.5

# Row 175 (future fixture - only has the odds columns filled in, no result yet)
$ws.Cells.Item(175, 1).Value = 173
$ws.Cells.Item(175, 2).Value = 6992698
$ws.Cells.Item(175, 3).Value = "Thailand Premier League"
$ws.Cells.Item(175, 4).Value = "Thailand Premier League"
$ws.Cells.Item(175, 5).Value = 45384.33333333334
$ws.Cells.Item(175, 6).Value = "Prachuap FC"
$ws.Cells.Item(175, 7).Value = "Trat FC"
$ws.Cells.Item(175, 8).ClearContents()
$ws.Cells.Item(175, 9).ClearContents()
$ws.Cells.Item(175, 10).ClearContents()
$ws.Cells.Item(175, 11).Value = 1.8
$ws.Cells.Item(175, 12).Value = 3.6
$ws.Cells.Item(175, 13).Value = 3.6
$ws.Cells.Item(175, 14).Value = 1.8
$ws.Cells.Item(175, 15).Value = 3.6
$ws.Cells.Item(175, 16).Value = 3.6
$ws.Cells.Item(175, 17).Value = -0.5
$ws.Cells.Item(175, 18).Value = 1.825
$ws.Cells.Item(175, 19).Value = 1.975
$ws.Cells.Item(175, 20).Value = 2.75
$ws.Cells.Item(175, 21).Value = 1.85
$ws.Cells.Item(175, 22).Value = 1.95
$ws.Cells.Item(175, 23).Value = 0
$ws.Cells.Item(175, 24).Value = 0
$ws.Cells.Item(175, 25).Value = 0
$ws.Cells.Item(175, 26).Value = 0
$ws.Cells.Item(175, 27).Value = 0
$ws.Cells.Item(175, 28).ClearContents()
$ws.Cells.Item(175, 29).ClearContents()
